$d = $word.ActiveDocument

# --- Step 1: remove the existing _GoBack bookmark (it will be re-added
#     at the end of the new second paragraph further down). ---
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# --- Step 2: fix up the end of paragraph 1. ---
# Before:  " của Nhất nghệ. Cá" (run) + "c bạn có thể tải full bộ tài liệu web php miễn phí tốt nhất." (run)
# After :  " của Nhất nghệ. Các bạn có thể tải full bộ tài liệu web php miễn phí tốt nhất." (single run)
$para1 = $d.Paragraphs(1).Range
$tailStart = $para1.Start + 147
$tailEnd = $para1.Start + 225
$tail = $d.Range($tailStart, $tailEnd)
$tail.Text = " của Nhất nghệ. Các bạn có thể tải full bộ tài liệu web php miễn phí tốt nhất."

# --- Step 3: append the two new paragraphs after paragraph 1. ---
$p1 = $d.Paragraphs(1).Range
$p1.InsertParagraphAfter()

$p2 = $d.Paragraphs(2).Range
$p2xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve">Tài liệu php chứa nhiều nội dung liên quan tới kỹ thuật lập trình website , quản trị , các kiến thức về lập trình web php từ </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>cơ bản đến nâng cao</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve">  hoàn toàn miển phí.</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
"@
$p2.InsertXML($p2xml)

$p2after = $d.Paragraphs(2).Range
$p2after.InsertParagraphAfter()

$p3 = $d.Paragraphs(3).Range
$p3xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
  </w:pPr>
</w:p>
"@
$p3.InsertXML($p3xml)

Write-Output ("ParaCount=" + $d.Paragraphs.Count)
